$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 571.89795
$ws.Range("I135").Value = 488.4565
$ws.Range("J135").Value = 1851.3334
$ws.Range("K135").Value = 4396.1085
$ws.Range("L135").Value = 16662.0006
$ws.Range("M135").Value = -1861.1085
$ws.Range("N135").Value = -21732.0006
$ws.Range("H138").Value = 1368.5892
$ws.Range("I138").Value = 1044.4131
$ws.Range("K138").Value = 3133.2393
$ws.Range("M138").Value = 2006.7607
$ws.Range("H141").Value = 1830.2128
$ws.Range("I141").Value = 1657.6136
$ws.Range("J141").Value = 4361.6665
$ws.Range("K141").Value = 4972.8408
$ws.Range("L141").Value = 13084.9995
$ws.Range("M141").Value = 207.1592000000001
$ws.Range("N141").Value = -23444.9995

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 15707
$ws.Range("J42").Value = 17266.666
$ws.Range("L42").Value = 17266.666
$ws.Range("N42").Value = -18238.666
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 2043.1904
$ws.Range("I61").Value = 1048.3334
$ws.Range("J61").Value = 3038.0476
$ws.Range("K61").Value = 1048.3334
$ws.Range("L61").Value = 3038.0476
$ws.Range("M61").Value = -836.3334
$ws.Range("N61").Value = -3462.0476
$ws.Range("H74").Value = 771.19446
$ws.Range("I74").Value = 763.76666
$ws.Range("K74").Value = 763.76666
$ws.Range("M74").Value = 110.23334
$ws.Range("H77").Value = 771.19446
$ws.Range("I77").Value = 763.76666
$ws.Range("K77").Value = 3818.8333
$ws.Range("M77").Value = 549.1666999999998
$ws.Range("H136").Value = 2043.1904
$ws.Range("I136").Value = 1048.3334
$ws.Range("J136").Value = 3038.0476
$ws.Range("K136").Value = 3145.0002
$ws.Range("L136").Value = 9114.1428
$ws.Range("M136").Value = -595.0002
$ws.Range("N136").Value = -14214.1428

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 42375.035
$ws.Range("I86").Value = 53010.09
$ws.Range("K86").Value = 53010.09
$ws.Range("M86").Value = -51887.09
$ws.Range("H89").Value = 42375.035
$ws.Range("I89").Value = 53010.09
$ws.Range("K89").Value = 265050.45
$ws.Range("M89").Value = -259434.45
$ws.Range("H107").Value = 66724176
$ws.Range("I107").Value = 76989370
$ws.Range("K107").Value = 76989370
$ws.Range("M107").Value = -76987450

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27146.91
$ws.Range("I31").Value = 1269.4073
$ws.Range("K31").Value = 1269.4073
$ws.Range("M31").Value = -974.4073000000001
$ws.Range("H34").Value = 27146.91
$ws.Range("I34").Value = 1269.4073
$ws.Range("K34").Value = 1269.4073
$ws.Range("M34").Value = -1067.4073
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H58").Value = 1006.09436
$ws.Range("I58").Value = 897.8571
$ws.Range("J58").Value = 2332
$ws.Range("K58").Value = 897.8571
$ws.Range("L58").Value = 2332
$ws.Range("M58").Value = -694.8571
$ws.Range("N58").Value = -2738
$ws.Range("H99").Value = 21416.285
$ws.Range("J99").Value = 41347.668
$ws.Range("L99").Value = 41347.668
$ws.Range("N99").Value = -44343.668
$ws.Range("H126").Value = 21416.285
$ws.Range("J126").Value = 41347.668
$ws.Range("L126").Value = 124043.004
$ws.Range("N126").Value = -128983.004
$ws.Range("H132").Value = 46879056
$ws.Range("I132").Value = 37040990
$ws.Range("K132").Value = 111122970
$ws.Range("M132").Value = -111120440
$ws.Range("H134").Value = 981.78
$ws.Range("I134").Value = 838.1539
$ws.Range("J134").Value = 1491
$ws.Range("K134").Value = 2514.4617
$ws.Range("L134").Value = 4473
$ws.Range("M134").Value = 20.53830000000016
$ws.Range("N134").Value = -9543
$ws.Range("H136").Value = 1006.09436
$ws.Range("I136").Value = 897.8571
$ws.Range("J136").Value = 2332
$ws.Range("K136").Value = 2693.5713
$ws.Range("L136").Value = 6996
$ws.Range("M136").Value = -143.5712999999996
$ws.Range("N136").Value = -12096

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1833.7778
$ws.Range("J51").Value = 2480
$ws.Range("L51").Value = 7440
$ws.Range("N51").Value = -8360
$ws.Range("H64").Value = 2525.3333
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 2961.1428
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 8883.4284
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -9423.4284
$ws.Range("H67").Value = 2525.3333
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 2961.1428
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 8883.4284
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -10755.4284
$ws.Range("H98").Value = 83683.836
$ws.Range("J98").Value = 200205.2
$ws.Range("L98").Value = 600615.6000000001
$ws.Range("N98").Value = -603611.6000000001
$ws.Range("H132").Value = 2019.091
$ws.Range("I132").Value = 1176.25
$ws.Range("J132").Value = 2500.7144
$ws.Range("K132").Value = 10586.25
$ws.Range("L132").Value = 22506.4296
$ws.Range("M132").Value = -8056.25
$ws.Range("N132").Value = -27566.4296

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58376.51
$ws.Range("I70").Value = 103727
$ws.Range("J70").Value = 5023
$ws.Range("K70").Value = 103727
$ws.Range("L70").Value = 5023
$ws.Range("M70").Value = -103457
$ws.Range("N70").Value = -5563
$ws.Range("H73").Value = 58376.51
$ws.Range("I73").Value = 103727
$ws.Range("J73").Value = 5023
$ws.Range("K73").Value = 103727
$ws.Range("L73").Value = 5023
$ws.Range("M73").Value = -102791
$ws.Range("N73").Value = -6895
$ws.Range("H102").Value = 2790.6667
$ws.Range("I102").Value = 2393.4285
$ws.Range("J102").Value = 3346.8
$ws.Range("K102").Value = 2393.4285
$ws.Range("L102").Value = 3346.8
$ws.Range("M102").Value = -771.4285
$ws.Range("N102").Value = -6590.8
$ws.Range("H132").Value = 8252.333
$ws.Range("I132").Value = 19000
$ws.Range("J132").Value = 2878.5
$ws.Range("K132").Value = 57000
$ws.Range("L132").Value = 8635.5
$ws.Range("M132").Value = -54470
$ws.Range("N132").Value = -13695.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3817.7058
$ws.Range("I7").Value = 2660
$ws.Range("J7").Value = 4300.0835
$ws.Range("K7").Value = 2660
$ws.Range("L7").Value = 4300.0835
$ws.Range("M7").Value = -2548
$ws.Range("N7").Value = -4524.0835
$ws.Range("H40").Value = 45385.434
$ws.Range("I40").Value = 143809.42
$ws.Range("J40").Value = 2324.9375
$ws.Range("K40").Value = 143809.42
$ws.Range("L40").Value = 2324.9375
$ws.Range("M40").Value = -143673.42
$ws.Range("N40").Value = -2596.9375
$ws.Range("H126").Value = 3817.7058
$ws.Range("I126").Value = 2660
$ws.Range("J126").Value = 4300.0835
$ws.Range("K126").Value = 7980
$ws.Range("L126").Value = 12900.2505
$ws.Range("M126").Value = -5510
$ws.Range("N126").Value = -17840.2505
$ws.Range("H132").Value = 2854.5757
$ws.Range("I132").Value = 2983.742
$ws.Range("J132").Value = 852.5
$ws.Range("K132").Value = 8951.226
$ws.Range("L132").Value = 2557.5
$ws.Range("M132").Value = -6421.226000000001
$ws.Range("N132").Value = -7617.5
$ws.Range("H135").Value = 35897.5
$ws.Range("J135").Value = 35897.5
$ws.Range("L135").Value = 35897.5
$ws.Range("N135").Value = -46037.5
$ws.Range("H136").Value = 1104.6428
$ws.Range("I136").Value = 950.3929
$ws.Range("J136").Value = 1413.1428
$ws.Range("K136").Value = 2851.1787
$ws.Range("L136").Value = 4239.428400000001
$ws.Range("M136").Value = -301.1787000000004
$ws.Range("N136").Value = -9339.4284

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2105.6667
$ws.Range("I126").Value = 2107.6
$ws.Range("J126").Value = 2103.25
$ws.Range("K126").Value = 6322.799999999999
$ws.Range("L126").Value = 6309.75
$ws.Range("M126").Value = -3852.799999999999
$ws.Range("N126").Value = -11249.75
$ws.Range("H132").Value = 2647.6
$ws.Range("I132").Value = 2645.6562
$ws.Range("J132").Value = 2668.3333
$ws.Range("K132").Value = 7936.9686
$ws.Range("L132").Value = 8004.999899999999
$ws.Range("M132").Value = -5406.9686
$ws.Range("N132").Value = -13064.9999
$ws.Range("H136").Value = 654.5686
$ws.Range("I136").Value = 378.55814
$ws.Range("K136").Value = 1135.67442
$ws.Range("M136").Value = 1414.32558
